$d = $word.ActiveDocument

# The roster table (table 2) lists team members. The edit:
#  - removes Youssef Reda Mokhtar entirely
#  - shifts everyone below him up by one slot (renumbering the leading
#    "N- " index in each name cell)
#  - the last row (previously "6- Rimaz Abd El-Rahman Emam") keeps its own
#    ID/Section/Group values; it is simply renumbered to "5-"
#  - the row that used to belong to Habiba Mohamed Fawzy is dropped from
#    the table since everyone shifted up and there is now one fewer row

function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $start = $cell.Range.Start
    $end = $cell.Range.End
    # Building a fresh Document.Range (rather than re-using the TableCell's
    # Range object directly) ensures the whole cell's contents -- even when
    # they span multiple runs -- gets replaced instead of only the first run.
    $d.Range($start, $end).Text = $newText
}

$t = $d.Tables.Item(2)

# Row 3 used to be Youssef's row; it becomes Abdulrahman's info (was row 4)
Set-CellText $t 3 1 "2- Abdulrahman Mohamed Mahmoud"
Set-CellText $t 3 2 "20220901"
Set-CellText $t 3 3 "3"
Set-CellText $t 3 4 "A"

# Row 4 used to be Abdulrahman's row; it becomes Menna's info (was row 5)
Set-CellText $t 4 1 "3- Menna Khaled Gamal"
Set-CellText $t 4 2 "20221166"
Set-CellText $t 4 3 "1"
Set-CellText $t 4 4 "A"

# Row 5 used to be Menna's row; it becomes Habiba's info (was row 6)
Set-CellText $t 5 1 "4- Habiba Mohamed Fawzy"
Set-CellText $t 5 2 "20220107"
Set-CellText $t 5 3 "1"
Set-CellText $t 5 4 "A"

# Row 6 used to be Habiba's row; it is removed entirely since Rimaz's row
# (previously row 7) now takes its place, just renumbered below.
$t.Rows.Item(6).Delete()

# Row 6 (now Rimaz's row, previously row 7) only needs its leading number
# changed from "6-" to "5-"; its ID/Section/Group stay the same.
Set-CellText $t 6 1 "5- Rimaz Abd El-Rahman Emam"
